$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.795.38'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '3.060.07'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.34'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.52'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = '3.058.49'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +3.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("E11").Value = '  -2.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.480'
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  +2.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.33'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").Value = '3.562.01'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").Value = '63.898.69'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").Value = '3.065.73'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.78'
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '486.48'
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").Value = '  +3.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.690'
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.69'
$ws.Range("E23").Value = '  +9.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.52'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.61'
$ws.Range("E25").Value = '  +2.38%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.18'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.50'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("E32").Value = '  +1.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.57'
$ws.Range("E33").Value = '  +3.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.69'
$ws.Range("E34").Value = '  +2.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.24'
$ws.Range("E35").Value = '  +2.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.83'
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("E37").Value = '  +1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '444.51'
$ws.Range("E38").Value = '  -4.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0814'
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("D40").Value = '3.049.27'
$ws.Range("E40").Value = '  +2.91%  '
$ws.Range("E41").Value = '  -5.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.33'
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("E43").Value = '  +2.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.275'
$ws.Range("E44").Value = '  +7.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '27.99'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("E46").Value = '  +6.11%  '
$ws.Range("E48").Value = '  +1.92%  '
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.83'
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.13'
$ws.Range("E51").Value = '  +3.88%  '
